# Updates cryptos list D/E columns per latest scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells look like plain numbers but must stay text
# (e.g. trailing zeros such as "562.60" must survive), so force
# the Text number format before writing the value - otherwise Excel
# would silently reinterpret them as numeric.
$textCells = @("D5","D6","D8","D11","D12","D14","D17","D20","D21","D24","D25","D30","D31","D33","D35","D38","D39","D41","D45","D47","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "59.216.72"
$ws.Range("E2").Value = "  +3.43%  "
$ws.Range("D3").Value = "2.999.17"
$ws.Range("E3").Value = "  +3.47%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "562.60"
$ws.Range("E5").Value = "  +2.60%  "
$ws.Range("D6").Value = "139.66"
$ws.Range("E6").Value = "  +10.90%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.519"
$ws.Range("E8").Value = "  +2.44%  "
$ws.Range("D9").Value = "2.986.46"
$ws.Range("E9").Value = "  +3.21%  "
$ws.Range("E10").Value = "  +7.58%  "
$ws.Range("D11").Value = "5.14"
$ws.Range("E11").Value = "  +9.18%  "
$ws.Range("D12").Value = "0.455"
$ws.Range("E12").Value = "  +5.09%  "
$ws.Range("E13").Value = "  +7.61%  "
$ws.Range("D14").Value = "33.76"
$ws.Range("E14").Value = "  +5.20%  "
$ws.Range("E15").Value = "  +2.74%  "
$ws.Range("D16").Value = "3.493.38"
$ws.Range("E16").Value = "  +3.54%  "
$ws.Range("D17").Value = "7.08"
$ws.Range("E17").Value = "  +8.62%  "
$ws.Range("D18").Value = "2.995.81"
$ws.Range("E18").Value = "  +3.82%  "
$ws.Range("D19").Value = "59.203.15"
$ws.Range("E19").Value = "  +3.54%  "
$ws.Range("D20").Value = "429.17"
$ws.Range("E20").Value = "  +6.15%  "
$ws.Range("D21").Value = "13.62"
$ws.Range("E21").Value = "  +6.42%  "
$ws.Range("E22").Value = "  +7.05%  "
$ws.Range("E23").Value = "  +4.33%  "
$ws.Range("D24").Value = "13.40"
$ws.Range("E24").Value = "  +5.87%  "
$ws.Range("D25").Value = "80.80"
$ws.Range("E25").Value = "  +3.98%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("E28").Value = "  +12.79%  "
$ws.Range("E29").Value = "  +3.51%  "
$ws.Range("D30").Value = "7.76"
$ws.Range("E30").Value = "  +7.16%  "
$ws.Range("D31").Value = "25.79"
$ws.Range("E31").Value = "  +4.86%  "
$ws.Range("E32").Value = "  +3.07%  "
$ws.Range("D33").Value = "0.0990"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("D34").Value = "0.0₃0777"
$ws.Range("E34").Value = "  +21.55%  "
$ws.Range("D35").Value = "0.995"
$ws.Range("E35").Value = "  +9.40%  "
$ws.Range("E36").Value = "  +7.82%  "
$ws.Range("E37").Value = "  +2.83%  "
$ws.Range("D38").Value = "49.32"
$ws.Range("E38").Value = "  +2.89%  "
$ws.Range("D39").Value = "8.65"
$ws.Range("E39").Value = "  +5.78%  "
$ws.Range("E40").Value = "  +11.47%  "
$ws.Range("D41").Value = "403.51"
$ws.Range("E41").Value = "  +10.94%  "
$ws.Range("D42").Value = "2.772.96"
$ws.Range("E42").Value = "  +6.44%  "
$ws.Range("E43").Value = "  +4.58%  "
$ws.Range("E44").Value = "  +1.81%  "
$ws.Range("D45").Value = "0.254"
$ws.Range("E45").Value = "  +11.24%  "
$ws.Range("D47").Value = "122.85"
$ws.Range("E47").Value = "  +2.29%  "
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("E49").Value = "  +3.92%  "
$ws.Range("D50").Value = "33.64"
$ws.Range("E50").Value = "  +23.12%  "
$ws.Range("E51").Value = "  +4.12%  "
